# Purchase Request workbook update:
# - Item 1 (row 14): Qty changes from 5 to 20, Description changes from "keyboard" to "Diesel"
# - Item 2 (row 15): removed entirely (Item No., Qty, UOM, Description cleared)
# - Active selection moves to A15:K15 (the now-empty row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 14 values
$ws.Range("B14").Value = 20
$ws.Range("E14").Value = "Diesel"

# Clear out row 15 (second line item) entirely
$ws.Range("A15:K15").ClearContents() | Out-Null

# Move the selection to the now-empty row, matching the saved view state
$ws.Range("A15:K15").Select() | Out-Null
